# RPA datasets push 2024-07-27
# - Adds a new NH-lead IPO row (엔에이치스팩31호) to the report, subscribed
#   2024-07-16 / paid-in 2024-07-19 / listed 2024-07-26.
# - The regenerated export also re-ordered two pre-existing rows relative to
#   the previous run (씨어스테크놀로지 stays put; the 한국제15호스팩 /
#   에스오에스랩 pair that follows it swaps places).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param(
        [int]$Row,
        [object[]]$Values,
        [bool[]]$AsText
    )
    for ($i = 0; $i -lt $Values.Length; $i++) {
        $col = $i + 1
        $val = $Values[$i]
        if ($AsText[$i]) {
            # Leading apostrophe forces text storage so date-shaped strings
            # (e.g. "2024-07-16") aren't auto-converted to date serials.
            $ws.Cells.Item($Row, $col).Value = "'" + $val
        } else {
            $ws.Cells.Item($Row, $col).Value = $val
        }
    }
    $ws.Range($ws.Cells.Item($Row, 1), $ws.Cells.Item($Row, $Values.Length)).ClearFormats()
}

# Column layout: A=인수기관 B=청약일 C=회사명 D=대표주관회사 E=인수회사
#                F=납입일 G=상장일 H=공모금액(백만원) I=공모주수 J=공모가
#                K=인수수수료 L=인수비율
$textFlags = @($true, $true, $true, $true, $true, $true, $true, $false, $false, $false, $false, $false)

# Insert the new row between the LS (이베스트스팩6호) row and the NH
# (시프트업) row, shifting every row below it down by one.
$ws.Rows.Item(7).Insert()

Set-RowValues -Row 7 -Values @(
    'NH', '2024-07-16', '엔에이치스팩31호', 'NH', 'NH', '2024-07-19', '2024-07-26',
    12000, 6000000, 2000, 0, 100
) -AsText $textFlags

# After the insert, the old rows 22/23 now sit at 23/24 — put the
# regenerated report's order back: 한국제15호스팩 first, then 에스오에스랩.
Set-RowValues -Row 23 -Values @(
    '한국', '2024-06-17', '한국제15호스팩', '한국', '한국', '2024-06-20', '2024-06-26',
    12500, 6250000, 2000, 0, 100
) -AsText $textFlags

Set-RowValues -Row 24 -Values @(
    '한국', '2024-06-14', '에스오에스랩', '한국', '한국, BNK', '2024-06-19', '2024-06-25',
    21850, 2000000, 11500, 0, 95
) -AsText $textFlags

Write-Host "RPA datasets push 2024-07-27 applied"
